$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- new data (previously belonging to the "Vedticka" record)
$ws.Range("A2").Value = 111487426
$ws.Range("B2").Value = 89369
$ws.Range("E2").Value = 5447
$ws.Range("F2").Value = "Vedticka"
$ws.Range("G2").Value = "Fuscoporia viticola"
$ws.Range("H2").Value = "(Schwein.) Murrill"
$ws.Range("Q2").Value = 626133.5793112689
$ws.Range("R2").Value = 6893051.461214696

# I2 becomes a present-but-empty text cell (quote-prefix trick forces an
# empty text value instead of clearing the cell outright), then strip the
# quote-prefix style back to Normal so no stray formatting is introduced.
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"

# AF2 is removed entirely in the edited file.
$ws.Range("AF2").ClearContents()

# Row 3 <- new data (previously belonging to the "Korallrot" record)
$ws.Range("A3").Value = 111487425
$ws.Range("B3").Value = 96251
$ws.Range("E3").Value = 220093
$ws.Range("F3").Value = "Korallrot"
$ws.Range("G3").Value = "Corallorhiza trifida"
$ws.Range("H3").Value = "Châtel."
$ws.Range("Q3").Value = 626157.6942840694
$ws.Range("R3").Value = 6893095.882089161

# I3 becomes the text "30" (must stay text, not be reinterpreted as a
# number), using the quote-prefix trick then clearing the resulting style.
$ws.Range("I3").Value = "'30"
$ws.Range("I3").Style = "Normal"

# AF3 is newly added as a present-but-empty text cell.
$ws.Range("AF3").Value = "'"
$ws.Range("AF3").Style = "Normal"
